$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# --- Content fixes -------------------------------------------------------

# Currency label / value clean-up
$wsInput.Range("A6").Value = "currency"

# Product name: add hyphen after "822" (same value mirrored on the output sheet)
$newProductName = "822-RBI-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-Late Repayment"
$wsInput.Range("B1").Value  = $newProductName
$wsOutput.Range("B1").Value = $newProductName

$wsInput.Range("B6").Value = "US Dollar"

# Give the currency value cell the same "accounting" look used further down
# the sheet (green fill, default font) instead of the old gray label style.
$wsInput.Range("B6").Interior.Color = 5296274

# --- Active sheet / selection state --------------------------------------

# ProductLoanOutput is no longer the active tab; its stored selection moves
# to B1.
$wsOutput.Activate()
$wsOutput.Range("B1").Select()

# ProductLoanInput becomes the active tab, with A6:B6 selected.
$wsInput.Activate()
$wsInput.Range("A6:B6").Select()
